$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Goal (per commit): add a space between "if" and "($money)" in the
# sentence "Geen spaties bijvoorbeeld if($money) en geen if ( $money )",
# and move the (hidden) "_GoBack" bookmark from the end of the document
# to the position right after the newly inserted space.
# ------------------------------------------------------------------

# Step 1: remove the existing "_GoBack" bookmark (it will be re-created
# at the new cursor position below, mirroring real Word's behaviour of
# moving _GoBack to the last edit location).
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

# Step 2: find the exact split point between "if" and "($money)".
$d1 = $word.ActiveDocument
$findRange = $d1.Content
$found = $findRange.Find.Execute("if(`$money)", $true, $false, $false, $false, `
                                  $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'if(`$money)' in the document"
}
$splitPos = $findRange.Start + 2   # right between "if" and "($money)"

# Step 3: stage a single space character inside an isolated, plain run
# ("nee" - paragraph 6) that carries no rsid/extra attributes, so that
# touching it to force a run-split leaves no visible trace once the
# staged character is cut back out. Toggling Bold off/on is a reliable
# way to force the freshly inserted character into its own <w:r>
# (Word/this host won't silently re-merge runs with different direct
# formatting history), after which we cut that lone run.
$stagePara = $d1.Paragraphs(6)
$stagePos = $stagePara.Range.End - 1   # just before the paragraph mark, end of "nee"
$stage = $d1.Range($stagePos, $stagePos)
$stage.InsertAfter(" ")
$stageRun = $d1.Range($stagePos, $stagePos + 1)
$stageRun.Font.Bold = 1
$stageRun.Font.Bold = 0
$stageRun2 = $d1.Range($stagePos, $stagePos + 1)
$stageRun2.Cut()

# Step 4: paste the staged space at the split point. Because the stage
# insert+cut is a net-zero length change, $splitPos is still valid.
$d2 = $word.ActiveDocument
$target = $d2.Range($splitPos, $splitPos)
$target.Paste()

# Step 5: re-create the "_GoBack" bookmark right after the pasted space
# (i.e. right before "($money)"), matching its new cursor position.
$d3 = $word.ActiveDocument
$bmPos = $splitPos + 1
$bmRange = $d3.Range($bmPos, $bmPos)
$d3.Bookmarks.Add("_GoBack", $bmRange)
